$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D previously carried a distinct "Arial-only" font style (s="2").
# The update flattens it to the same style used by the rest of the sheet.
$ws.Range("D1:D5").Font.Name = "Arial"

# Row 3 ("dim"-level" metadata row): orden/siglas columns get real values
# instead of placeholder "null" / mismatched label.
$ws.Range("B3").Value = "iaest-measure:orden"
$ws.Range("D3").Value = "iaest-measure:siglas"

# Row 4 ("medida"-level metadata row): orden/siglas columns get "medida"
# instead of "null" / "dim".
$ws.Range("B4").Value = "medida"
$ws.Range("D4").Value = "medida"

# Row 5 (type/URI metadata row): orden/siglas columns get their proper
# XSD type labels instead of "null" / "skos:Concept".
$ws.Range("B5").Value = "xsd:int"
$ws.Range("D5").Value = "xsd:string"

# Row 6 (the stray "mapping-siglas.xlsx" note in D6) is removed entirely.
$ws.Rows.Item(6).Delete()
